$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.647.29"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.509.48"
$ws.Range("E3").Value = "  -0.07%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.42"
$ws.Range("E5").Value = "  -0.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.68"
$ws.Range("E6").Value = "  -0.26%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.509.99"
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.162"
$ws.Range("E10").Value = "  +1.07%  "
$ws.Range("E11").Value = "  -0.46%  "
$ws.Range("E12").Value = "  +3.92%  "
$ws.Range("E13").Value = "  +1.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.966.39"
$ws.Range("E14").Value = "  -0.75%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "69.576.80"
$ws.Range("E15").Value = "  +0.05%  "
$ws.Range("E16").Value = "  +1.47%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.89"
$ws.Range("E17").Value = "  -0.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.510.47"
$ws.Range("E18").Value = "  -0.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.27"
$ws.Range("E19").Value = "  -1.56%  "
$ws.Range("E20").Value = "  -3.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "350.08"
$ws.Range("E21").Value = "  -0.13%  "
$ws.Range("E22").Value = "  -1.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.94"
$ws.Range("E23").Value = "  -0.85%  "
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "70.14"
$ws.Range("E25").Value = "  +1.95%  "
$ws.Range("E26").Value = "  -1.49%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.83"
$ws.Range("E27").Value = "  -1.87%  "
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0891"
$ws.Range("E30").Value = "  -1.30%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.85"
$ws.Range("E31").Value = "  -0.91%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "461.02"
$ws.Range("E32").Value = "  -3.54%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.22"
$ws.Range("E33").Value = "  -5.02%  "
$ws.Range("E34").Value = "  -1.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "0.117"
$ws.Range("E36").Value = "  +0.44%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").Value = "158.72"
$ws.Range("E37").Value = "  +2.38%  "
$ws.Range("E38").Value = "  +0.80%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.48"
$ws.Range("E39").Value = "  -0.46%  "
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("E41").Value = "  -0.22%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.69"
$ws.Range("E42").Value = "  -1.53%  "
$ws.Range("E43").Value = "  -0.58%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "38.18"
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.22"
$ws.Range("E45").Value = "  -4.75%  "
$ws.Range("E46").Value = "  -7.99%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "142.44"
$ws.Range("E47").Value = "  -1.22%  "
$ws.Range("E48").Value = "  -2.10%  "
$ws.Range("E49").Value = "  -2.46%  "
$ws.Range("E50").Value = "  +0.44%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.581"
$ws.Range("E51").Value = "  -1.12%  "
